$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple numeric corrections scattered throughout the existing data ---
$ws.Range("C95").Value = 5
$ws.Range("C100").Value = 3
$ws.Range("C1102").Value = 29
$ws.Range("C1145").Value = 48
$ws.Range("C1156").Value = 42
$ws.Range("C1199").Value = 56
$ws.Range("C1260").Value = 14
$ws.Range("C1279").Value = 20
$ws.Range("C1342").Value = 7
$ws.Range("C1369").Value = 11
$ws.Range("C1372").Value = 4
$ws.Range("C1375").Value = 21
$ws.Range("C1387").Value = 5
$ws.Range("C1394").Value = 11
$ws.Range("C1404").Value = 4
$ws.Range("C1406").Value = 11

# --- Rework of the trailing rows: the last few days of data were revised
#     and several new days (2021-02-13 .. 2021-02-16) were appended. ---

# Row 1408 (2021-02-13 / 44240): age group changes from 60-69 to 40-49
$ws.Range("B1408").Value = "40-49"

# Row 1409 (2021-02-13 / 44240): age group changes from 70-79 to 60-69, count 1 -> 3
$ws.Range("B1409").Value = "60-69"
$ws.Range("C1409").Value = 3

# Row 1410 (2021-02-13 / 44240): age group changes from 80+ to 70-79
$ws.Range("B1410").Value = "70-79"

# Row 1411 previously held 2021-02-14 data; now becomes the new 80+ row for 2021-02-13
$ws.Range("A1411").Value = 44240
$ws.Range("B1411").Value = "80+"
$ws.Range("C1411").Value = 7

# New number format needs to be applied to freshly-created date cells (A1412:A1421)
# so they match the existing "YYYY-MM-DD HH:MM:SS" date style used by column A.
$ws.Range("A1412:A1421").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1412: 2021-02-14
$ws.Range("A1412").Value = 44241
$ws.Range("B1412").Value = "40-49"
$ws.Range("C1412").Value = 2

# Row 1413: 2021-02-14
$ws.Range("A1413").Value = 44241
$ws.Range("B1413").Value = "50-59"
$ws.Range("C1413").Value = 1

# Row 1414: 2021-02-14
$ws.Range("A1414").Value = 44241
$ws.Range("B1414").Value = "60-69"
$ws.Range("C1414").Value = 3

# Row 1415: 2021-02-14
$ws.Range("A1415").Value = 44241
$ws.Range("B1415").Value = "70-79"
$ws.Range("C1415").Value = 1

# Row 1416: 2021-02-14
$ws.Range("A1416").Value = 44241
$ws.Range("B1416").Value = "80+"
$ws.Range("C1416").Value = 3

# Row 1417: 2021-02-15
$ws.Range("A1417").Value = 44242
$ws.Range("B1417").Value = "50-59"
$ws.Range("C1417").Value = 1

# Row 1418: 2021-02-15
$ws.Range("A1418").Value = 44242
$ws.Range("B1418").Value = "60-69"
$ws.Range("C1418").Value = 4

# Row 1419: 2021-02-15
$ws.Range("A1419").Value = 44242
$ws.Range("B1419").Value = "70-79"
$ws.Range("C1419").Value = 4

# Row 1420: 2021-02-15
$ws.Range("A1420").Value = 44242
$ws.Range("B1420").Value = "80+"
$ws.Range("C1420").Value = 10

# Row 1421: 2021-02-16
$ws.Range("A1421").Value = 44243
$ws.Range("B1421").Value = "70-79"
$ws.Range("C1421").Value = 2
